# Fixed update to excel issue
# 1) Rename the "Requested quantity" header on the Weekly/Monthly sheets
# 2) Add a new "PO Forecast" sheet with the forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header style (bold/bordered/centered) from the Weekly sheet's
# header row so the new header cells pick up the same shared style (s=1)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the date-formatted style (s=2) used on column A of the Weekly sheet
# down the whole of column A on the new sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A49").PasteSpecial(-4122)

$data = @(
    @(44934.99999999999, 0, -586.695198591931, 455.9768803723563),
    @(44997.99999999999, 0, -540.9631063160023, 546.5547262994863),
    @(45011.99999999999, 14, -488.2956572692235, 568.5673707267588),
    @(45039.99999999999, 47, -485.5767096281834, 576.7357747172656),
    @(45060.99999999999, 71, -448.6330532146201, 564.5938216762546),
    @(45067.99999999999, 80, -442.9017232445206, 662.5236741025965),
    @(45081.99999999999, 96, -394.8099705143273, 614.2716996153378),
    @(45088.99999999999, 104, -388.6443448174869, 666.8111214967297),
    @(45095.99999999999, 112, -378.2108335087352, 609.6904019992932),
    @(45102.99999999999, 120, -397.3945066658862, 642.435051592313),
    @(45109.99999999999, 128, -413.5154937886398, 631.900290285196),
    @(45116.99999999999, 137, -358.6709921381384, 696.9157204645103),
    @(45123.99999999999, 145, -384.9267116542648, 683.8181725194975),
    @(45130.99999999999, 153, -400.4758040761587, 652.6223604183159),
    @(45137.99999999999, 161, -376.6746522741856, 665.4114272852627),
    @(45144.99999999999, 169, -319.7121294208431, 679.560681437242),
    @(45193.99999999999, 226, -302.6171586517008, 752.0191849567319),
    @(45221.99999999999, 259, -245.50072533473, 790.5284664960062),
    @(45228.99999999999, 267, -233.2938404289193, 794.2424923538256),
    @(45235.99999999999, 275, -276.1101207099678, 812.9347856871266),
    @(45242.99999999999, 283, -212.3647908292633, 815.4306685594752),
    @(45270.99999999999, 316, -218.1144970907055, 810.8831528009782),
    @(45277.99999999999, 324, -213.3562911543891, 873.2623622489782),
    @(45298.99999999999, 348, -152.9082038806744, 855.794777641715),
    @(45305.99999999999, 356, -143.6192491119223, 908.0175001698223),
    @(45312.99999999999, 365, -168.0051741245328, 888.4032098500347),
    @(45319.99999999999, 373, -162.3862263231135, 902.3544897219801),
    @(45326.99999999999, 381, -159.5625413403811, 917.6352488397785),
    @(45333.99999999999, 389, -108.8092333691167, 901.0601266849558),
    @(45340.99999999999, 397, -118.6903855652987, 939.3169352149719),
    @(45347.99999999999, 405, -104.7885940764996, 936.1805119594604),
    @(45354.99999999999, 413, -72.54419365374233, 958.8187441900966),
    @(45361.99999999999, 422, -160.7412824236568, 947.9131859403129),
    @(45368.99999999999, 430, -104.5932763318241, 989.0379658346241),
    @(45375.99999999999, 438, -121.4513190313096, 954.0624195036984),
    @(45550.99999999999, 642, 124.7238161300927, 1155.87567947777),
    @(45557.99999999999, 650, 112.8341034210497, 1193.369463533929),
    @(45578.99999999999, 674, 139.1086217763086, 1180.6602455049),
    @(45585.99999999999, 682, 159.1140370025745, 1244.432481760005),
    @(45592.99999999999, 690, 103.3652645160806, 1221.398334577317),
    @(45599.99999999999, 699, 122.5572188947559, 1213.970035312552),
    @(45606.99999999999, 707, 184.383109678738, 1232.767079730363),
    @(45613.99999999999, 715, 201.7669132960641, 1268.492996910115),
    @(45620.99999999999, 723, 207.5053044235603, 1230.658943912468),
    @(45627.99999999999, 731, 192.1586333140868, 1251.621444864233),
    @(45634.99999999999, 739, 221.9605957428697, 1267.926202419286),
    @(45641.99999999999, 747, 250.5002618527631, 1272.502650665131),
    @(45648.99999999999, 756, 192.6518290615449, 1243.160321020464)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value2 = $row[0]
    $wsForecast.Cells.Item($r, 2).Value2 = $row[1]
    $wsForecast.Cells.Item($r, 3).Value2 = $row[2]
    $wsForecast.Cells.Item($r, 4).Value2 = $row[3]
    $r = $r + 1
}

# Restore the originally active sheet/selection
$wsWeekly.Activate()
